# Auto-generated: update LeveProfits market-price-derived columns (H-N)
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(97, 8).Value = 44723.2
$ws.Cells.Item(97, 9).Value = 950
$ws.Cells.Item(97, 10).Value = 55666.5
$ws.Cells.Item(97, 11).Value = 2850
$ws.Cells.Item(97, 12).Value = 166999.5
$ws.Cells.Item(97, 13).Value = -2354
$ws.Cells.Item(97, 14).Value = -167991.5
$ws.Cells.Item(132, 8).Value = 1532.081
$ws.Cells.Item(132, 9).Value = 1262.6333
$ws.Cells.Item(132, 11).Value = 3787.8999
$ws.Cells.Item(132, 13).Value = -1257.8999
$ws.Cells.Item(137, 8).Value = 1647.3478
$ws.Cells.Item(137, 9).Value = 1494.1538
$ws.Cells.Item(137, 11).Value = 4482.4614
$ws.Cells.Item(137, 13).Value = -1932.4614
$ws.Cells.Item(138, 8).Value = 5884825.5
$ws.Cells.Item(138, 9).Value = 1243.0714
$ws.Cells.Item(138, 10).Value = 8775007
$ws.Cells.Item(138, 11).Value = 3729.2142
$ws.Cells.Item(138, 12).Value = 26325021
$ws.Cells.Item(138, 13).Value = 1410.7858
$ws.Cells.Item(138, 14).Value = -26335301

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3031.0715
$ws.Cells.Item(2, 9).Value = 2353.5
$ws.Cells.Item(2, 10).Value = 4725
$ws.Cells.Item(2, 11).Value = 2353.5
$ws.Cells.Item(2, 12).Value = 4725
$ws.Cells.Item(2, 13).Value = -2240.5
$ws.Cells.Item(2, 14).Value = -4951
$ws.Cells.Item(32, 8).Value = 10343.643
$ws.Cells.Item(32, 9).Value = 6295.8857
$ws.Cells.Item(32, 11).Value = 6295.8857
$ws.Cells.Item(32, 13).Value = -6008.8857
$ws.Cells.Item(116, 8).Value = 3031.0715
$ws.Cells.Item(116, 9).Value = 2353.5
$ws.Cells.Item(116, 10).Value = 4725
$ws.Cells.Item(116, 11).Value = 2353.5
$ws.Cells.Item(116, 12).Value = 4725
$ws.Cells.Item(116, 13).Value = -59.5
$ws.Cells.Item(116, 14).Value = -9313
$ws.Cells.Item(122, 8).Value = 1276.0769
$ws.Cells.Item(122, 9).Value = 895.8889
$ws.Cells.Item(122, 10).Value = 2131.5
$ws.Cells.Item(122, 11).Value = 2687.6667
$ws.Cells.Item(122, 12).Value = 6394.5
$ws.Cells.Item(122, 13).Value = -237.6667000000002
$ws.Cells.Item(122, 14).Value = -11294.5
$ws.Cells.Item(132, 8).Value = 3023.725
$ws.Cells.Item(132, 9).Value = 2921.7693
$ws.Cells.Item(132, 10).Value = 7000
$ws.Cells.Item(132, 11).Value = 8765.3079
$ws.Cells.Item(132, 12).Value = 21000
$ws.Cells.Item(132, 13).Value = -6235.3079
$ws.Cells.Item(132, 14).Value = -26060

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3031.0715
$ws.Cells.Item(3, 9).Value = 2353.5
$ws.Cells.Item(3, 10).Value = 4725
$ws.Cells.Item(3, 11).Value = 2353.5
$ws.Cells.Item(3, 12).Value = 4725
$ws.Cells.Item(3, 13).Value = -2239.5
$ws.Cells.Item(3, 14).Value = -4953
$ws.Cells.Item(107, 8).Value = 1257.25
$ws.Cells.Item(107, 9).Value = 1151.1428
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 1151.1428
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 768.8571999999999
$ws.Cells.Item(107, 14).Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 254.17647
$ws.Cells.Item(7, 9).Value = 296.22223
$ws.Cells.Item(7, 10).Value = 206.875
$ws.Cells.Item(7, 11).Value = 296.22223
$ws.Cells.Item(7, 12).Value = 206.875
$ws.Cells.Item(7, 13).Value = -183.22223
$ws.Cells.Item(7, 14).Value = -432.875
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).ClearContents()
$ws.Cells.Item(33, 8).Value = 21945.889
$ws.Cells.Item(33, 9).Value = 3906
$ws.Cells.Item(33, 10).Value = 44495.75
$ws.Cells.Item(33, 11).Value = 3906
$ws.Cells.Item(33, 12).Value = 44495.75
$ws.Cells.Item(33, 13).Value = -3527
$ws.Cells.Item(33, 14).Value = -45253.75
$ws.Cells.Item(42, 8).Value = 5000
$ws.Cells.Item(42, 9).Value = 5000
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 5000
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -4407
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(55, 8).Value = 41650
$ws.Cells.Item(55, 9).Value = 29900
$ws.Cells.Item(55, 11).Value = 29900
$ws.Cells.Item(55, 13).Value = -29585
$ws.Cells.Item(56, 8).Value = 1166.6666
$ws.Cells.Item(56, 9).Value = 750
$ws.Cells.Item(56, 10).Value = 2000
$ws.Cells.Item(56, 11).Value = 750
$ws.Cells.Item(56, 12).Value = 2000
$ws.Cells.Item(56, 13).Value = 95
$ws.Cells.Item(56, 14).Value = -3690
$ws.Cells.Item(99, 8).Value = 3248.5
$ws.Cells.Item(99, 9).Value = 3248.5
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 3248.5
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -1750.5
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(105, 8).Value = 1165.6
$ws.Cells.Item(105, 10).Value = 1120.4445
$ws.Cells.Item(105, 12).Value = 1120.4445
$ws.Cells.Item(105, 14).Value = -4614.4445
$ws.Cells.Item(126, 8).Value = 3248.5
$ws.Cells.Item(126, 9).Value = 3248.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 9745.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -7275.5
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2709.6345
$ws.Cells.Item(132, 9).Value = 2618.02
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 7854.059999999999
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -5324.059999999999
$ws.Cells.Item(132, 14).Value = -20060
$ws.Cells.Item(134, 8).Value = 12662.47
$ws.Cells.Item(134, 9).Value = 4930.8047
$ws.Cells.Item(134, 10).Value = 52287.25
$ws.Cells.Item(134, 11).Value = 14792.4141
$ws.Cells.Item(134, 12).Value = 156861.75
$ws.Cells.Item(134, 13).Value = -12257.4141
$ws.Cells.Item(134, 14).Value = -161931.75
$ws.Cells.Item(141, 8).Value = 341448.88
$ws.Cells.Item(141, 10).Value = 341448.88
$ws.Cells.Item(141, 12).Value = 341448.88
$ws.Cells.Item(141, 14).Value = -351808.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 29417490
$ws.Cells.Item(4, 9).Value = 30282706
$ws.Cells.Item(4, 11).Value = 90848118
$ws.Cells.Item(4, 13).Value = -90848006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 5000
$ws.Cells.Item(5, 10).Value = 5000
$ws.Cells.Item(5, 12).Value = 5000
$ws.Cells.Item(5, 14).Value = -5224
$ws.Cells.Item(122, 8).Value = 2553.9
$ws.Cells.Item(122, 9).Value = 2215.158
$ws.Cells.Item(122, 10).Value = 8990
$ws.Cells.Item(122, 11).Value = 6645.474
$ws.Cells.Item(122, 12).Value = 26970
$ws.Cells.Item(122, 13).Value = -4195.474
$ws.Cells.Item(122, 14).Value = -31870
$ws.Cells.Item(132, 8).Value = 2414.818
$ws.Cells.Item(132, 10).Value = 6628.25
$ws.Cells.Item(132, 12).Value = 19884.75
$ws.Cells.Item(132, 14).Value = -24944.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 6681.3076
$ws.Cells.Item(13, 9).Value = 6317.1113
$ws.Cells.Item(13, 11).Value = 6317.1113
$ws.Cells.Item(13, 13).Value = -6177.1113
$ws.Cells.Item(61, 8).Value = 3884.5652
$ws.Cells.Item(61, 9).Value = 3746.95
$ws.Cells.Item(61, 10).Value = 4802
$ws.Cells.Item(61, 11).Value = 3746.95
$ws.Cells.Item(61, 12).Value = 4802
$ws.Cells.Item(61, 13).Value = -3544.95
$ws.Cells.Item(61, 14).Value = -5206
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 3884.5652
$ws.Cells.Item(113, 9).Value = 3746.95
$ws.Cells.Item(113, 10).Value = 4802
$ws.Cells.Item(113, 11).Value = 3746.95
$ws.Cells.Item(113, 12).Value = 4802
$ws.Cells.Item(113, 13).Value = -1576.95
$ws.Cells.Item(113, 14).Value = -9142

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 250000
$ws.Cells.Item(2, 9).Value = 250000
$ws.Cells.Item(2, 11).Value = 250000
$ws.Cells.Item(2, 13).Value = -249888
$ws.Cells.Item(95, 8).Value = 99999.25
$ws.Cells.Item(95, 10).Value = 99999.25
$ws.Cells.Item(95, 12).Value = 99999.25
$ws.Cells.Item(95, 14).Value = -105491.25
$ws.Cells.Item(100, 8).Value = 619.0454999999999
$ws.Cells.Item(100, 10).Value = 718.4
$ws.Cells.Item(100, 12).Value = 1436.8
$ws.Cells.Item(100, 14).Value = -2518.8
$ws.Cells.Item(129, 8).Value = 39669.332
$ws.Cells.Item(129, 10).Value = 39669.332
$ws.Cells.Item(129, 12).Value = 39669.332
$ws.Cells.Item(129, 14).Value = -49669.332
$ws.Cells.Item(130, 8).Value = 61713.5
$ws.Cells.Item(130, 10).Value = 61713.5
$ws.Cells.Item(130, 12).Value = 61713.5
$ws.Cells.Item(130, 14).Value = -71753.5
$ws.Cells.Item(132, 8).Value = 1611.7587
$ws.Cells.Item(132, 9).Value = 1508.9259
$ws.Cells.Item(132, 11).Value = 4526.7777
$ws.Cells.Item(132, 13).Value = -1996.7777
$ws.Cells.Item(138, 8).Value = 57166.668
$ws.Cells.Item(138, 10).Value = 57166.668
$ws.Cells.Item(138, 12).Value = 57166.668
$ws.Cells.Item(138, 14).Value = -67446.66800000001
